$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look like plain numbers,
# so Excel keeps them as text (matching original inlineStr string cells)
# instead of silently converting them to numeric values.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'

# Apply the updated values
$ws.Range('D2').Value = '37.367.69'
$ws.Range('E2').Value = '  +2.24%  '

$ws.Range('D3').Value = '2.066.14'
$ws.Range('E3').Value = '  +3.47%  '

$ws.Range('E4').Value = '  -0.01%  '

$ws.Range('D5').Value = '235.64'
$ws.Range('E5').Value = '  +0.01%  '

$ws.Range('D6').Value = '0.613'
$ws.Range('E6').Value = '  +2.19%  '

$ws.Range('D7').Value = '58.61'
$ws.Range('E7').Value = '  +6.90%  '

$ws.Range('E8').Value = '  -0.01%  '

$ws.Range('D9').Value = '0.383'
$ws.Range('E9').Value = '  +3.09%  '

$ws.Range('D10').Value = '58.99'
$ws.Range('E10').Value = '  +1.33%  '

$ws.Range('D11').Value = '0.0763'
$ws.Range('E11').Value = '  +1.95%  '

$ws.Range('E12').Value = '  +2.67%  '

$ws.Range('D13').Value = '2.369.87'
$ws.Range('E13').Value = '  +3.52%  '

$ws.Range('D14').Value = '14.61'
$ws.Range('E14').Value = '  +3.06%  '

$ws.Range('D15').Value = '21.41'
$ws.Range('E15').Value = '  +4.97%  '

$ws.Range('D16').Value = '0.778'
$ws.Range('E16').Value = '  +2.76%  '

$ws.Range('D17').Value = '5.20'
$ws.Range('E17').Value = '  +2.28%  '

$ws.Range('D18').Value = '2.046.24'
$ws.Range('E18').Value = '  +2.33%  '

$ws.Range('D19').Value = '37.560.18'
$ws.Range('E19').Value = '  +2.80%  '

$ws.Range('D20').Value = '6.13'
$ws.Range('E20').Value = '  +15.96%  '

$ws.Range('D21').Value = '70.47'
$ws.Range('E21').Value = '  +3.88%  '

$ws.Range('D22').Value = '0.0₃0814'
$ws.Range('E22').Value = '  +1.05%  '

$ws.Range('D23').Value = '227.95'
$ws.Range('E23').Value = '  +2.65%  '

$ws.Range('E24').Value = '  +0.03%  '

$ws.Range('D25').Value = '2.45'
$ws.Range('E25').Value = '  +1.41%  '

$ws.Range('D26').Value = '2.40'
$ws.Range('E26').Value = '  +1.13%  '

$ws.Range('D27').Value = '165.68'
$ws.Range('E27').Value = '  +2.17%  '

$ws.Range('E28').Value = '  +12.54%  '

$ws.Range('D29').Value = '8.88'
$ws.Range('E29').Value = '  +2.45%  '

$ws.Range('E30').Value = '  +0.74%  '

$ws.Range('D31').Value = '19.21'
$ws.Range('E31').Value = '  +1.67%  '

$ws.Range('E32').Value = '  +1.84%  '

$ws.Range('D33').Value = '4.52'
$ws.Range('E33').Value = '  +3.03%  '

$ws.Range('D34').Value = '0.0623'
$ws.Range('E34').Value = '  +3.07%  '

$ws.Range('E35').Value = '  +8.49%  '

$ws.Range('D36').Value = '4.55'
$ws.Range('E36').Value = '  +6.98%  '

$ws.Range('D37').Value = '3.39'
$ws.Range('E37').Value = '  +1.15%  '

$ws.Range('D38').Value = '0.999'
$ws.Range('E38').Value = '  -0.06%  '

$ws.Range('E39').Value = '  +1.33%  '

$ws.Range('D40').Value = '5.87'
$ws.Range('E40').Value = '  +4.05%  '

$ws.Range('D41').Value = '0.0983'
$ws.Range('E41').Value = '  +3.69%  '

$ws.Range('E42').Value = '  -1.32%  '

$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D43').Value = '4.36'
$ws.Range('E43').Value = '  +23.22%  '

$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '1.457.13'
$ws.Range('E44').Value = '  +0.29%  '

$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = '95.54'
$ws.Range('E45').Value = '  +7.28%  '

$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '1.17'
$ws.Range('E46').Value = '  +5.45%  '

$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').Value = '0.0212'
$ws.Range('E47').Value = '  +4.33%  '

$ws.Range('D48').Value = '15.87'
$ws.Range('E48').Value = '  +3.96%  '

$ws.Range('E49').Value = '  +3.80%  '

$ws.Range('E50').Value = '  +5.59%  '

$ws.Range('E51').Value = '  +1.92%  '

# Remove the temporary text-number format so no stray style is left behind,
# while keeping the values stored as text.
$ws.Range('D5').ClearFormats()
$ws.Range('D6').ClearFormats()
$ws.Range('D7').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D11').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D33').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
